# ANA-1050.xlsx -> adds a new "InsertServiceRequest" apex snippet for the
# "Closed" status in cell K3 (row 3), matching the new K1/K2 pair already
# present for the "Open" status. Mirrors test-case ANA-1053 creation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cell: K3 holds the apex snippet that inserts a Closed
# SVMXC__Service_Request__c record (sibling of K2's "Open" variant).
$ws.Range("K3").Value = "SVMXC__Service_Request__c SR_1 = new SVMXC__Service_Request__c(SVMXC__Status__c = 'Closed' );insert SR_1 ;"

# Row 3 now needs extra height so the long wrapped string is readable,
# matching the row-2 precedent (which already carries an explicit height).
$ws.Rows.Item(3).RowHeight = 48

# Move the cursor / view to the newly added cell.
$ws.Range("K3").Select()
